# feat: leader board topic-xp-stats
#
# Rewrites the exercise data grid (rows 2-6, columns A-K) on Sheet1 with the
# new leader-board / topic-xp-stats sample rows, clears cells that are no
# longer used, and moves the active selection to K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SELECT_IMAGE
$ws.Range("A2").Value = "SELECT_IMAGE"
$ws.Range("B2").Value = "Look and match: 'hi'"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "g3-u1-hi.png"
$ws.Range("I2").Value = "g3-u1-bye.png"
$ws.Range("J2").Value = "g1-u2-car.png"
$ws.Range("K2").Value = "g1-u2-car.png"

# Row 3 - MULTIPLE_CHOICE
$ws.Range("A3").Value = "MULTIPLE_CHOICE"
$ws.Range("B3").Value = "_i"
$ws.Range("C3").Value = "g3-u1-hi.png"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "Hi"
$ws.Range("G3").Value = "Xin chào"
$ws.Range("H3").Value = "h"
$ws.Range("I3").Value = "b"
$ws.Range("J3").Value = "c"
$ws.Range("K3").Value = "d"

# Row 4 - LISTENING
$ws.Range("A4").Value = "LISTENING"
$ws.Range("B4").Value = "Listen and choose"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "g3-u1-hi.mp3"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "Hi"
$ws.Range("G4").Value = "Xin chào"
$ws.Range("H4").Value = "Hi"
$ws.Range("I4").Value = "Bye"
$ws.Range("J4").Value = "Goobye"
$ws.Range("K4").Value = "Hello"

# Row 5 - PRONUNCIATION
$ws.Range("A5").Value = "PRONUNCIATION"
$ws.Range("B5").Value = "Listen and pronunciation"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "g3-u1-hi.mp3"
$ws.Range("E5").Value = "Hi"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""

# Row 6 - MATCHING
$ws.Range("A6").Value = "MATCHING"
$ws.Range("B6").Value = "Matching"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "Hi|xin chào"
$ws.Range("I6").Value = "Hi Image|g3-u1-hi.png"
$ws.Range("J6").Value = "g1-u2-car.png|g1-u2-car.png"
$ws.Range("K6").Value = "Bill|Hello"

# Move the active selection to K5, matching the saved view state.
$ws.Range("K5").Select()
